# Update the Role for the second user (row 2) from "Student" to "Professor"
# and move the active selection to F2, matching the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Professor"

$ws.Range("F2").Select()
